$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row height + cell content updates for rows 7-18 ---
$ws.Rows.Item(7).RowHeight = 91
$ws.Range("A7").Value = "subir archivo de datos"
$ws.Range("B7").Value = "/archivo/?version"
$ws.Range("C7").Value = "sube un conjunto de datos al servidor especificando la version del archivo de configuracion"
$ws.Range("D7").Value = "post"
$ws.Range("E7").Value = "version"
$ws.Range("F7").Value = "{`narchivo: {objeto_tipo_datos}`n}"
$ws.Range("G7").Value = "{`n status: ok, `nmensaje: “creado con exito”`n}"
$ws.Range("H7").Value = "revisar el archivo de especificaciones para ver exactamente como esta constituido este archivo"

$ws.Rows.Item(8).RowHeight = 91
$ws.Range("A8").Value = "obtener archivo de datos"
$ws.Range("B8").Value = "/archivo/?version"
$ws.Range("C8").Value = "obtiene el archivo de datos de la version especificada desde el servidor"
$ws.Range("D8").Value = "get"
$ws.Range("E8").Value = "version"
$ws.Range("G8").Value = "regresa_archivo"
$ws.Range("H8").Value = "revisar el archivo de especificaciones para ver exactamente como esta constituido este archivo"

$ws.Rows.Item(9).RowHeight = 68.65
$ws.Range("A9").Value = "guardar (cambiar nombre a archivo) archivo de datos"
$ws.Range("B9").Value = "/archivo/guardar"
$ws.Range("C9").Value = "guarda el archivo de datos con un nombre especificado por el usuario"
$ws.Range("D9").Value = "post"
$ws.Range("F9").Value = "{`narchivo_anterior: “nombre_archivo”,`narchivo_nuevo: “nombre_archivo”`n}"
$ws.Range("G9").Value = "{`n status: ok, `nmensaje: “guardado con exito”`n}"

$ws.Rows.Item(10).RowHeight = 91.5
$ws.Range("A10").Value = "crear instancia (nueva fila de datos)"
$ws.Range("B10").Value = "/datos/?version/"
$ws.Range("C10").Value = "crea un nuevo objeto en el archivo"
$ws.Range("D10").Value = "post"
$ws.Range("E10").Value = "version"
$ws.Range("F10").Value = "{`nobjeto: objeto`n}"
$ws.Range("G10").Value = "{`n status: ok, `nmensaje: “creado con exito”`n}"
$ws.Range("H10").Value = "revisar el archivo de especificaciones para ver exactamente como esta constituido este archivo"

$ws.Rows.Item(11).RowHeight = 91.5
$ws.Range("A11").Value = "modificar instancia (modificar algun valor o valores de la fila)"
$ws.Range("B11").Value = "/archivo/?version"
$ws.Range("C11").Value = "modifica un objeto en el archivo"
$ws.Range("D11").Value = "patch"
$ws.Range("E11").Value = "version"
$ws.Range("F11").Value = "{`nid: id,`nobjeto: objeto`n}"
$ws.Range("G11").Value = "{`n status: ok, `nmensaje: “modificado con exito”`n}"
$ws.Range("H11").Value = "revisar el archivo de especificaciones para ver exactamente como esta constituido este archivo"

$ws.Rows.Item(12).RowHeight = 91
$ws.Range("A12").Value = "eliminar instancia (eliminar algun objeto (fila) del archivo)"
$ws.Range("B12").Value = "/archivo/?version/?id"
$ws.Range("C12").Value = "elimina un objeto completo del archivo"
$ws.Range("D12").Value = "delete"
$ws.Range("E12").Value = "version,id"
$ws.Range("G12").Value = "{`n status: ok, `nmensaje: “eliminado con exito”`n}"
$ws.Range("H12").Value = "revisar el archivo de especificaciones para ver exactamente como esta constituido el id"

$ws.Rows.Item(13).RowHeight = 57.75
$ws.Range("A13").Value = "agregar atributo (columna completa)"
$ws.Range("B13").Value = "/archivo/?version/atributo"
$ws.Range("C13").Value = "crea una columna completa en todos los datos (con el valor especificado como sin valor)"
$ws.Range("D13").Value = "post"
$ws.Range("E13").Value = "version"
$ws.Range("F13").Value = "{`nnombre_atributo: “nombre”`n}"
$ws.Range("G13").Value = "{`n status: ok, `nmensaje: “creado con exito”`n}"

$ws.Rows.Item(14).RowHeight = 57.75
$ws.Range("A14").Value = "eliminar atributo (columna  completa)"
$ws.Range("B14").Value = "/archivo/?version/atributo/?nombre"
$ws.Range("C14").Value = "elimina una columna completa (en todos los objetos del archivo)"
$ws.Range("D14").Value = "delete"
$ws.Range("E14").Value = "version, nombre"
$ws.Range("G14").Value = "{`n status: ok, `nmensaje: “eliminado con exito”`n}"

$ws.Rows.Item(15).RowHeight = 57.75
$ws.Range("A15").Value = "obtener bases de datos disponibles"
$ws.Range("B15").Value = "/base-de-datos"
$ws.Range("C15").Value = "obtiene el nombre de las bases de datos disponibles"
$ws.Range("D15").Value = "get"
$ws.Range("G15").Value = "{`n status: ok, `nbases: [“nombre1”, “nombre2”, “etc”]`n}"

$ws.Rows.Item(16).RowHeight = 57.75
$ws.Range("A16").Value = "obtener tablas de las bases de datos"
$ws.Range("B16").Value = "/base-de-datos/?nombre/tablas"
$ws.Range("C16").Value = "obtiene el nombre de las tablas de la base de datos especificada "
$ws.Range("D16").Value = "get"
$ws.Range("E16").Value = "nombre"
$ws.Range("G16").Value = "{`n status: ok, `ntablas: [“nombre1”, “nombre2”, “etc”]`n}"

$ws.Rows.Item(17).RowHeight = 57.75
$ws.Range("A17").Value = "obtener atributos de la tabla"
$ws.Range("B17").Value = "/base-de-datos/?nombre/tablas/?nombre-tabla"
$ws.Range("C17").Value = "obtiene la lista de los atributos de la tabla especificada de la base de datos especificada"
$ws.Range("D17").Value = "get"
$ws.Range("E17").Value = "nombre, nombre-tabla"
$ws.Range("G17").Value = "{`n status: ok, `natributos: [“nombre1”, “nombre2”, “etc”]`n}"

$ws.Rows.Item(18).RowHeight = 57.75
$ws.Range("A18").Value = "obtener datos de la tabla"
$ws.Range("B18").Value = "/base-de-datos/?nombre/?nombre-tabla"
$ws.Range("C18").Value = "obtiene los datos de  la tabla especificada"
$ws.Range("D18").Value = "get"
$ws.Range("E18").Value = "nombre, nombre-tabla"
$ws.Range("G18").Value = "{`n status: ok, `ndatos: objeto_tipo_datos`n}"

# --- Placeholder "incomplete" cells with underline style (H9, F14) ---
$ws.Range("H9").WrapText = $true
$ws.Range("H9").Font.Underline = $true
$ws.Range("F14").WrapText = $true
$ws.Range("F14").Font.Underline = $true

# --- Column width adjustments (col B widened, col E added width) ---
$ws.Columns.Item(2).ColumnWidth = 30.166666666666668
$ws.Columns.Item(5).ColumnWidth = 14.5

# --- Selection / scroll position ---
$excel.ActiveWindow.ScrollRow = 14
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A19").Select()
